$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "69.438.01"
$ws.Range("E2").Value = "  +1.85%  "

# Row 3
$ws.Range("D3").Value = "3.901.62"
$ws.Range("E3").Value = "  +0.53%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
Set-TextValue "D5" "529.61"
$ws.Range("E5").Value = "  +9.78%  "

# Row 6
Set-TextValue "D6" "145.27"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("E7").Value = "  -1.23%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
Set-TextValue "D9" "0.720"
$ws.Range("E9").Value = "  -2.61%  "

# Row 10
$ws.Range("E10").Value = "  -1.94%  "

# Row 11
$ws.Range("E11").Value = "  -4.01%  "

# Row 12
Set-TextValue "D12" "42.28"
$ws.Range("E12").Value = "  -1.59%  "

# Row 13
$ws.Range("D13").Value = "4.521.07"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14
Set-TextValue "D14" "10.28"
$ws.Range("E14").Value = "  -1.65%  "

# Row 15
$ws.Range("D15").Value = "3.895.54"
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
Set-TextValue "D16" "13.99"
$ws.Range("E16").Value = "  -1.31%  "

# Row 17
$ws.Range("E17").Value = "  -1.18%  "

# Row 18
$ws.Range("E18").Value = "  +7.16%  "

# Row 19
Set-TextValue "D19" "19.85"
$ws.Range("E19").Value = "  -0.30%  "

# Row 20
$ws.Range("D20").Value = "69.329.29"
$ws.Range("E20").Value = "  +1.65%  "

# Row 21
Set-TextValue "D21" "425.27"

# Row 22
Set-TextValue "D22" "3.40"
$ws.Range("E22").Value = "  -4.62%  "

# Row 23
$ws.Range("E23").Value = "  -3.72%  "

# Row 24
Set-TextValue "D24" "88.22"
$ws.Range("E24").Value = "  -0.54%  "

# Row 25
Set-TextValue "D25" "4.06"
$ws.Range("E25").Value = "  +10.56%  "

# Row 26
Set-TextValue "D26" "11.41"
$ws.Range("E26").Value = "  -7.32%  "

# Row 27
Set-TextValue "D27" "10.60"
$ws.Range("E27").Value = "  -3.23%  "

# Row 28
Set-TextValue "D28" "36.41"
$ws.Range("E28").Value = "  -1.74%  "

# Row 29
Set-TextValue "D29" "691.44"
$ws.Range("E29").Value = "  -3.85%  "

# Row 30
Set-TextValue "D30" "13.21"
$ws.Range("E30").Value = "  -1.73%  "

# Row 31
$ws.Range("E31").Value = "  -2.64%  "

# Row 32
Set-TextValue "D32" "2.84"
$ws.Range("E32").Value = "  -2.69%  "

# Row 33
Set-TextValue "D33" "68.98"
$ws.Range("E33").Value = "  +11.71%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("E34").Value = "  -0.54%  "

# Row 35
$ws.Range("E35").Value = "  +9.02%  "

# Row 36
$ws.Range("E36").Value = "  -1.46%  "

# Row 37
Set-TextValue "D37" "40.10"
$ws.Range("E37").Value = "  -1.38%  "

# Row 38
$ws.Range("E38").Value = "  +2.68%  "

# Row 39
$ws.Range("E39").Value = "  -0.02%  "

# Row 40
$ws.Range("E40").Value = "  -0.20%  "

# Row 41
Set-TextValue "D41" "3.24"
$ws.Range("E41").Value = "  +5.70%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0482"
$ws.Range("E42").Value = "  -2.60%  "

# Row 43
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D43" "3.19"
$ws.Range("E43").Value = "  +8.42%  "

# Row 44
Set-TextValue "D44" "2.80"
$ws.Range("E44").Value = "  -5.97%  "

# Row 45
$ws.Range("E45").Value = "  +1.59%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D46" "0.141"
$ws.Range("E46").Value = "  -0.75%  "

# Row 47
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D47" "0.000286"
$ws.Range("E47").Value = "  +15.15%  "

# Row 48
Set-TextValue "D48" "3.00"
$ws.Range("E48").Value = "  +7.05%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0347"
$ws.Range("E49").Value = "  -2.25%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "146.08"
$ws.Range("E50").Value = "  +1.28%  "

# Row 51
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.744.27"
$ws.Range("E51").Value = "  +15.07%  "
